$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 9 trailing stale data rows (the sheet shrinks from 290 to 281 data rows).
# Deleting them shifts the existing blank row / "Filtros aplicados" note row up
# so they keep their place immediately after the data block.
$ws.Range("A283:C291").EntireRow.Delete()

# Final "Conta / Nome / Saldo" data, one pipe-delimited row per line:
#   account number | account holder name | balance
$rows = @(
"004214592|MERG|236549.41",
"004565108|GUSTAVO|153046.59",
"004213929|RODOLFO|130000.12",
"004222784|RAFAEL|129988.31",
"004582648|ANUAR|110803.31",
"004556974|KELLY|109472.42",
"004499920|FABIANO|97011.13",
"004363260|LARISSA|84936.84",
"004453157|JOSE|77685.88",
"004461070|EDUARDO|72290.24",
"004212581|MARIA|67263.08",
"005233971|LEIVANIO|63595.66",
"004935287|ODILON|51109.31",
"004398253|EULER|50000",
"004612043|YURI|50000",
"004508526|CASSIO|42835.62",
"004550605|REJANE|38896.93",
"004973881|ISABELLA|34288.29",
"004805333|ROSANA|33629.79",
"004376145|LUCYENE|33022.32",
"004872395|ALESSANDRA|31515.95",
"004550750|THEO|28955.69",
"004484207|FLAVIA|23121.33",
"004313254|GUSTAVO|18656.83",
"005366255|RAPHAELA|18000",
"004237325|RICARDO|17000",
"005064129|THIAGO|15844.88",
"002064834|RAFAELA|14820.88",
"004565146|GUSTAVO|14530.99",
"005324840|PEDRO|13987.81",
"004454365|RAFAEL|13701.88",
"004454491|MARCO|11855.11",
"005007114|STELLA|10235.26",
"005076418|LEONARDO|9682.02",
"005040864|ANDRE|8028.72",
"004165515|MAURO|5535.68",
"005341184|BRENO|5193.98",
"004361159|HFR|5087.64",
"004221638|CAROLINE|4605.84",
"005063749|NATALIA|4500",
"001368670|THIAGO|3299.19",
"004224815|GUILHERME|3067.85",
"004572740|PAULO|3000",
"002277249|DANILO|2239.1",
"004884046|WILSON|1328.21",
"004435987|MARCO|1073.38",
"004488571|CARLOS|1000",
"004472431|LUIS|860.7",
"004749680|OTAVIANO|855.95",
"004381328|JOAO|851.22",
"004563252|FERNANDO|776.3",
"004381180|HFR|743.31",
"004547722|MARCIA|700",
"004216504|WANDER|681.62",
"004387250|MONICA|648.53",
"004813088|JULIANA|644.25",
"004646727|RENATA|637.96",
"004644524|PAULO|632.63",
"004804125|EDUARDO|608.35",
"005079458|JONAS|594.89",
"004452912|BRUNO|562.08",
"004212409|RAFAEL|554.54",
"004643746|MARIO|552.33",
"005142592|ALBERTO|551.27",
"004468717|HELOISA|470.9",
"004210959|ANA|463.15",
"005003629|ANDRE|402.57",
"004119016|HEMAT|399.92",
"005061124|BRUNO|384.02",
"004216401|SUELY|378.82",
"004862746|CESAR|375.85",
"004334062|MERG|350.28",
"004289402|LARISSA|349.86",
"004907688|HEITOR|312.22",
"004870976|HFR|309.15",
"004436055|MARCO|299.85",
"000626491|FELIPE|280.02",
"004754056|BRUNO|263.65",
"004556853|MARCEL|260.27",
"004751154|CATARINE|250.72",
"004479287|ANA|241.71",
"004472386|GABRIEL|230.83",
"004526450|MSD|205.84",
"005009947|VERANICE|199.46",
"004452946|OTAVIO|196.51",
"004767746|ISABELE|195.1",
"004475395|DAVID|185.02",
"002828327|RENAN|161.13",
"005022526|ALEXANDRE|147.18",
"004204500|EDWARD|140.07",
"004482090|CEZAR|127.5",
"004381359|PEDRO|124.54",
"004505474|RICARDO|118.95",
"004403674|LILIAN|116.45",
"004364200|BLOCO|107.23",
"005198093|ANA|105.51",
"004514241|ANDRE|104.73",
"002694089|VITOR|99.85",
"004908680|ELENE|99.31",
"004339183|JALISON|95.69",
"004334158|LEONE|93.87",
"004536602|TATIANY|93.49",
"004350197|GISELA|91.94",
"004228456|FLASH|91.65",
"005206566|LEVI|91",
"004239387|LUIZ|89.82",
"004243043|SUELI|89.79",
"004207278|CESAR|88.52",
"004459461|INTERLAGOS|86.7",
"004212132|JOAO|86.38",
"003249855|MARINA|86.02",
"005305965|SIDMAR|85.35",
"004360431|CARLOS|85.25",
"004473942|DAIANNE|84.23",
"005135281|RAFAEL|84.16",
"004950193|GABRIEL|83.5",
"004927044|CINTIA|81.49",
"004259659|BENTO|81.21",
"004318604|RENAN|80.51",
"004451996|ADRIANO|80.36",
"004267976|E3|79.84",
"001294033|VIVIANE|79.82",
"005009922|ANA|79.02",
"005338054|ELAINE|78.53",
"005070742|JUSCELINO|78.29",
"004643880|GABRIEL|78.23",
"004470679|RODOLFO|77.51",
"004477812|DIEGO|77.1",
"005274028|RAFAEL|76.97",
"004335144|EDMUNDO|76.27",
"004214460|MARIA|75.73",
"004229526|EDUARDA|74.99",
"004431591|MARIO|73.59",
"004399832|EULER|73.2",
"004713953|ALESSANDRA|73.04",
"005255637|PATRICIA|72.9",
"004915243|POLO|72.62",
"004460491|PEDRO|68.19",
"004332103|JOSE|68.01",
"003836362|ISABELLA|67.74",
"005133039|PAULO|66.51",
"004397124|MURYLO|65.64",
"004855596|MARIANA|64.36",
"004335251|EDMUNDO|62.39",
"004756968|DANIELY|61.45",
"004426743|GABRIELLE|59.28",
"004200433|BENTO|58.44",
"003894173|ANDREA|56.53",
"004517506|LUIZ|55.87",
"004215217|CAROLINA|55.66",
"004321092|DANIEL|53.54",
"004265173|JULIA|53.34",
"005032151|ANA|52.9",
"004479734|RODRIGO|52.24",
"004848279|LEONARDO|52.2",
"004749928|NADY|52.11",
"004259650|BENTO|51.65",
"004971448|CLOVIS|51.39",
"004115403|HEBERT|50.87",
"004517080|TATIANA|50.15",
"004208447|LEILA|50",
"004809902|PEDRO|49.49",
"004332747|LOHRAN|49.19",
"004413537|CLAUDIA|49.16",
"004388077|WLADMIR|48.42",
"004329229|GABRIEL|48.38",
"004392159|RODRIGO|48.16",
"004374943|LEONARDO|46.78",
"005216881|RENAN|46.76",
"003497496|ELISANDRA|46.23",
"004277637|LARA|46.2",
"005000460|MARIANA|46.09",
"004398174|DANIELE|45.63",
"004335031|EDMUNDO|45.3",
"004491730|DENISE|45.11",
"004497875|HENRIQUE|44.73",
"001731007|GUILHERME|44.59",
"004581652|CINCO|44.13",
"004805133|PATRICIA|41.48",
"004958578|ASSAKO|40.7",
"002697806|CLAUDIA|40.23",
"004998717|GIOVANE|40.08",
"003115072|VICTOR|39.35",
"004497825|PRISCILLA|38.79",
"004290978|LARISSA|38.7",
"003435941|HEITOR|38.51",
"004238164|DANIELA|38.3",
"005000656|LUCIA|35.88",
"004340984|RENATA|34",
"005044389|CLAUDIA|33.81",
"004994036|BALTASAR|33.73",
"004691225|ANNA|33.64",
"004376853|ALBERTO|30.98",
"004381095|THIAGO|30.49",
"004972070|MARIA|29.91",
"004481463|MARA|29.77",
"004472760|SANDRA|28.66",
"004486497|ELENA|28.63",
"004374891|RODRIGO|28.04",
"004755204|FABIANA|27.99",
"005018038|ELAINE|27.91",
"005105172|VALDIVINO|26.42",
"002738211|MARGARETH|26.3",
"001000882|AYRTON|25.52",
"004643153|CARLA|25.37",
"005245032|ROSA|25.18",
"004404724|LEANDRO|24.14",
"004630773|NABOR|24.1",
"005173958|VENIA|23.27",
"004452507|DANIELA|22.74",
"004278033|DAISY|21.59",
"005203796|MARCIA|20.88",
"004751770|DILSON|20.76",
"004214604|MARIA|20.72",
"004920447|MARILIA|20.67",
"004975924|SERGIO|20.11",
"004493324|DANIEL|20.01",
"004752615|LUZIMAR|19.9",
"004242237|MARIAH|19.06",
"004204255|AMADO|18.77",
"004368994|CRISTINA|18.56",
"004381194|ALINNE|17.71",
"005143579|GABRIEL|16.18",
"005169333|EDUARDO|16.12",
"004268684|PATRICIA|15.41",
"004512434|CAIO|15.39",
"004213943|ELISA|15.13",
"004224011|THOMAS|14.92",
"004422594|WANDIR|14.67",
"000827730|LUCIANA|13.29",
"004754920|LUIS|12.38",
"004936634|LEONARDO|10.81",
"004752461|SERGIO|10.77",
"001719494|LUIS|10.24",
"005103059|WALQUIRIA|10",
"004216298|FLORDELIZ|9.74",
"004527606|MARCIA|9.52",
"005002457|ROSANGELA|9.5",
"004921978|ELAINE|8.08",
"004381415|JOAO|7.95",
"004693631|NELY|7.36",
"005228239|DEBORA|7.11",
"004530494|ROSANGELA|6.94",
"004854496|JOSE|6.64",
"004332207|IRACY|5.87",
"004448501|JOAO|5.55",
"005142624|RODRIGO|4.75",
"004207658|ROBERTO|4.54",
"004504449|KELMA|3.65",
"004848927|ULDARICO|3.62",
"005142661|SABRINA|3.6",
"004425261|THAYSA|3.59",
"004575632|ADELE|3.46",
"005046805|RICARDO|3.23",
"004240292|MARCO|2.66",
"004886366|RENATO|1.57",
"004520100|ALEXANDRE|1.46",
"004264780|MARCELO|1.42",
"004308815|ZELI|1.25",
"004459875|HELVECIO|1.08",
"004452597|LARA|0.89",
"004223502|BRUNA|0.78",
"000834301|MARCUS|0.68",
"004380749|ELEUSE|0.6",
"005165116|ANA|0.51",
"004453302|ISABELLA|0.39",
"004924605|ESTER|0.3",
"004587511|CARLOS|0.24",
"004806286|VERA|0.19",
"004371857|NAZARETH|0.18",
"004357159|JOAO|0.15",
"004320840|NATALIA|0.14",
"004466350|RAQUEL|0.11",
"004806244|CARLA|0.1",
"005047946|GABRIEL|0.09",
"004589311|CLARICE|0.06",
"004473718|LUCAS|0.02",
"004850070|RENATO|0.02",
"002878817|GUILHERME|0.01",
"004400000|VILMA|0.01",
"005009992|ALINE|-4823.92"
)

# Account numbers must stay text so the leading zeros survive the round-trip.
$ws.Range("A2:A282").NumberFormat = "@"

$r = 2
foreach ($line in $rows) {
    $parts = $line.Split("|")
    $ws.Cells.Item($r, 1).Value = $parts[0]
    $ws.Cells.Item($r, 2).Value = $parts[1]
    $ws.Cells.Item($r, 3).Value = [double]$parts[2]
    $r = $r + 1
}
